$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns A, Q, R, S, AW, AX between row 2 and row 3.
$cols = @("A", "Q", "R", "S", "AW", "AX")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"

    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2

    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
